# Latest build of the game.
# Update the Capstone Hours log:
#  - Row 82's "Hours" cell changes from the shared text "2+" to the literal
#    number 2, and the old "2+" entry is reused as a brand new, more
#    descriptive task description that is moved down onto row 83.
#  - Two new rows (83 and 84) are appended with new task descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 82: "Hours" column becomes a plain numeric value instead of the
# old free-text "2+" entry.
$ws.Range("B82").Value = 2

# Row 83: new task entry that reuses/repurposes the old "2+" text slot.
$ws.Range("A83").Value = "Writing, Recording, putting in game new End Scene audio"
$ws.Range("B83").Value = 2
$ws.Range("C83").Value = 41950
$ws.Range("C82").Copy()
$ws.Range("C83").PasteSpecial(-4122)

# Row 84: another new task entry.
$ws.Range("A84").Value = "Playtesting, Finding bugs, fixing fatal bugs that made game non-playable"
$ws.Range("B84").Value = 7
$ws.Range("C84").Value = 41950
$ws.Range("C82").Copy()
$ws.Range("C84").PasteSpecial(-4122)

# Scroll/selection bookkeeping to mirror the author's view position.
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
$ws.Range("E79").Select()
